$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet 1 / index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4667
$ws1.Range("F6").Value = 3128
$ws1.Range("F9").Value = 273
$ws1.Range("F10").Value = 636
$ws1.Range("F12").Value = 536
$ws1.Range("F13").Value = 389
$ws1.Range("F14").Value = 137
$ws1.Range("F16").Value = 1350
$ws1.Range("F18").Value = 1623
$ws1.Range("F19").Value = 12
$ws1.Range("F26").Value = 52
$ws1.Range("F27").Value = 106
$ws1.Range("F32").Value = 3883
$ws1.Range("F36").Value = 983
$ws1.Range("F38").Value = 1857

# Sheet "演出" (sheet 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 47

# Sheet "全部类型" (sheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4667
$ws4.Range("F6").Value = 3128
$ws4.Range("F9").Value = 273
$ws4.Range("F10").Value = 636
$ws4.Range("F12").Value = 536
$ws4.Range("F14").Value = 389
$ws4.Range("F15").Value = 137
$ws4.Range("F17").Value = 1350
$ws4.Range("F19").Value = 1623
$ws4.Range("F20").Value = 12
$ws4.Range("F27").Value = 52
$ws4.Range("F28").Value = 106
$ws4.Range("F33").Value = 3883
$ws4.Range("F34").Value = 47
$ws4.Range("F38").Value = 983
$ws4.Range("F40").Value = 1857
